$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to keep a text/string type even when the value
    # looks like a number (e.g. "0.9979"), while leaving the cell's
    # resolved style index unchanged (reset back to the default "Normal"
    # style right after the write so no stray style is left behind).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "25.909.79"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.744.72"
$ws.Range("E3").Value = "  -0.35%  "
Set-TextValue "D4" "0.9979"
$ws.Range("E4").Value = "  -0.38%  "
Set-TextValue "D5" "229.83"
$ws.Range("E5").Value = "  -2.98%  "
Set-TextValue "D6" "0.9970"
$ws.Range("E6").Value = "  -0.40%  "
Set-TextValue "D7" "0.5155"
$ws.Range("E7").Value = "  +2.05%  "
Set-TextValue "D8" "0.2819"
$ws.Range("E8").Value = "  +8.06%  "
Set-TextValue "D9" "39.37"
$ws.Range("E9").Value = "  -2.36%  "
Set-TextValue "D10" "0.06114"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "1.740.88"
$ws.Range("E11").Value = "  -0.61%  "
Set-TextValue "D12" "0.06988"
Set-TextValue "D13" "15.37"
$ws.Range("E13").Value = "  +0.03%  "
Set-TextValue "D14" "0.6376"
$ws.Range("E14").Value = "  +5.76%  "
Set-TextValue "D15" "4.512"
$ws.Range("E15").Value = "  +1.52%  "
Set-TextValue "D16" "76.66"
$ws.Range("E16").Value = "  -2.03%  "
Set-TextValue "D17" "0.9979"
$ws.Range("E17").Value = "  -0.33%  "
Set-TextValue "D18" "0.9978"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "25.908.68"
$ws.Range("E19").Value = "  +0.46%  "
Set-TextValue "D20" "11.50"
$ws.Range("E20").Value = "  -1.07%  "
Set-TextValue "D21" "0.000006620"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "1.964.78"
$ws.Range("E22").Value = "  -0.65%  "
Set-TextValue "D23" "4.118"
$ws.Range("E23").Value = "  +1.70%  "
Set-TextValue "D24" "8.550"
$ws.Range("E24").Value = "  +4.56%  "
Set-TextValue "D25" "5.138"
$ws.Range("E25").Value = "  -0.26%  "
Set-TextValue "D26" "139.75"
$ws.Range("E26").Value = "  +1.31%  "
Set-TextValue "D27" "1.511"
$ws.Range("E27").Value = "  +3.18%  "
Set-TextValue "D28" "15.08"
$ws.Range("E28").Value = "  +0.07%  "
Set-TextValue "D29" "1.822"
$ws.Range("E29").Value = "  +1.35%  "
Set-TextValue "D30" "103.01"
$ws.Range("E30").Value = "  +1.06%  "
Set-TextValue "D31" "0.08306"
$ws.Range("E31").Value = "  +0.56%  "
Set-TextValue "D32" "3.631"
$ws.Range("E32").Value = "  -1.92%  "
Set-TextValue "D33" "3.422"
$ws.Range("E33").Value = "  +1.13%  "
Set-TextValue "D34" "0.04395"
$ws.Range("E34").Value = "  +0.96%  "
Set-TextValue "D35" "2.614"
$ws.Range("E35").Value = "  -1.39%  "
Set-TextValue "D36" "0.9778"
$ws.Range("E36").Value = "  -2.13%  "
Set-TextValue "D37" "0.6079"
$ws.Range("E37").Value = "  +1.71%  "
Set-TextValue "D38" "2.676"
$ws.Range("E38").Value = "  -0.50%  "
Set-TextValue "D39" "0.01565"
$ws.Range("E39").Value = "  +1.34%  "
Set-TextValue "D40" "1.928"
$ws.Range("E40").Value = "  -1.28%  "
Set-TextValue "D41" "0.9967"
$ws.Range("E41").Value = "  -0.41%  "
Set-TextValue "D42" "100.54"
$ws.Range("E42").Value = "  -2.37%  "
Set-TextValue "D43" "0.3843"
$ws.Range("E43").Value = "  +1.20%  "
Set-TextValue "D44" "0.7251"
$ws.Range("E44").Value = "  -3.18%  "
Set-TextValue "D45" "4.938"
$ws.Range("E45").Value = "  +1.98%  "
Set-TextValue "D46" "0.05452"
$ws.Range("E46").Value = "  -0.57%  "
Set-TextValue "D47" "6.358"
$ws.Range("E47").Value = "  +7.57%  "
Set-TextValue "D48" "0.1109"
$ws.Range("E48").Value = "  +3.21%  "
Set-TextValue "D49" "52.60"
$ws.Range("E49").Value = "  +1.47%  "
Set-TextValue "D50" "29.81"
$ws.Range("E50").Value = "  -0.95%  "
Set-TextValue "D51" "7.511"
$ws.Range("E51").Value = "  +1.12%  "
